$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("List1")
try {
  $s = $wb.Styles.Add("MyStyle")
  $ws1.Range("AA1").Style = "MyStyle"
  $ws1.Range("AA1").Value = "stylework"
} catch {
  $ws1.Range("AA1").Value = "styleerr: $($_.Exception.Message)"
}
